$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 19986246
$ws.Range("I33").Value = 24982694
$ws.Range("K33").Value = 24982694
$ws.Range("M33").Value = -24982465

$ws.Range("H40").Value = 1454.5714
$ws.Range("I40").Value = 1336.4
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 1336.4
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -1161.4
$ws.Range("N40").Value = -2100

$ws.Range("H41").Value = 352.7619
$ws.Range("I41").Value = 393.86667
$ws.Range("J41").Value = 250
$ws.Range("K41").Value = 393.86667
$ws.Range("L41").Value = 250
$ws.Range("M41").Value = 46.13333
$ws.Range("N41").Value = -1130

$ws.Range("H134").Value = 66961
$ws.Range("J134").Value = 66961
$ws.Range("L134").Value = 66961
$ws.Range("N134").Value = -77101

$ws.Range("H135").Value = 46875756
$ws.Range("I135").Value = 25000444
$ws.Range("K135").Value = 225003996
$ws.Range("M135").Value = -225001461

$ws.Range("H137").Value = 727047.4
$ws.Range("I137").Value = 2247.9614
$ws.Range("J137").Value = 1669286.6
$ws.Range("K137").Value = 6743.8842
$ws.Range("L137").Value = 5007859.800000001
$ws.Range("M137").Value = -4193.8842
$ws.Range("N137").Value = -5012959.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19077.45
$ws.Range("I32").Value = 23717.217
$ws.Range("J32").Value = 3832.5
$ws.Range("K32").Value = 23717.217
$ws.Range("L32").Value = 3832.5
$ws.Range("M32").Value = -23430.217
$ws.Range("N32").Value = -4406.5

$ws.Range("H44").Value = 65000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 65000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 65000
$ws.Range("N44").Value = -65976
$ws.Range("M44").ClearContents()

$ws.Range("H45").Value = 1491.6364
$ws.Range("I45").Value = 1412.762
$ws.Range("J45").Value = 1746.4615
$ws.Range("K45").Value = 1412.762
$ws.Range("L45").Value = 1746.4615
$ws.Range("M45").Value = -1035.762
$ws.Range("N45").Value = -2500.4615

$ws.Range("H55").Value = 70000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 70000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 70000
$ws.Range("N55").Value = -70630
$ws.Range("M55").ClearContents()

$ws.Range("H80").Value = 40105
$ws.Range("J80").Value = 40105
$ws.Range("L80").Value = 40105
$ws.Range("N80").Value = -42101

$ws.Range("H83").Value = 40105
$ws.Range("J83").Value = 40105
$ws.Range("L83").Value = 120315
$ws.Range("N83").Value = -130299

$ws.Range("H132").Value = 2499.6758
$ws.Range("I132").Value = 2402.6667
$ws.Range("J132").Value = 3300
$ws.Range("K132").Value = 7208.000100000001
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -4678.000100000001
$ws.Range("N132").Value = -14960

$ws.Range("H133").Value = 52820.332
$ws.Range("J133").Value = 52820.332
$ws.Range("L133").Value = 52820.332
$ws.Range("N133").Value = -57880.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 245.57895
$ws.Range("I80").Value = 123.14286
$ws.Range("J80").Value = 317
$ws.Range("K80").Value = 123.14286
$ws.Range("L80").Value = 317
$ws.Range("M80").Value = 874.85714
$ws.Range("N80").Value = -2313

$ws.Range("H82").Value = 14723.706
$ws.Range("I82").Value = 4074
$ws.Range("J82").Value = 40283
$ws.Range("K82").Value = 4074
$ws.Range("L82").Value = 40283
$ws.Range("M82").Value = -3691
$ws.Range("N82").Value = -41049

$ws.Range("H83").Value = 245.57895
$ws.Range("I83").Value = 123.14286
$ws.Range("J83").Value = 317
$ws.Range("K83").Value = 615.7143
$ws.Range("L83").Value = 1585
$ws.Range("M83").Value = 4376.2857
$ws.Range("N83").Value = -11569

$ws.Range("H85").Value = 14723.706
$ws.Range("I85").Value = 4074
$ws.Range("J85").Value = 40283
$ws.Range("K85").Value = 4074
$ws.Range("L85").Value = 40283
$ws.Range("M85").Value = -2748
$ws.Range("N85").Value = -42935

$ws.Range("H107").Value = 3024.6428
$ws.Range("I107").Value = 2921.2
$ws.Range("K107").Value = 2921.2
$ws.Range("M107").Value = -1001.2

$ws.Range("H134").Value = 70484.734
$ws.Range("I134").Value = 3135.7
$ws.Range("J134").Value = 205182.8
$ws.Range("K134").Value = 9407.099999999999
$ws.Range("L134").Value = 615548.3999999999
$ws.Range("M134").Value = -6872.099999999999
$ws.Range("N134").Value = -620618.3999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2212.2222
$ws.Range("I16").Value = 1985
$ws.Range("J16").Value = 2666.6667
$ws.Range("K16").Value = 1985
$ws.Range("L16").Value = 2666.6667
$ws.Range("M16").Value = -1698
$ws.Range("N16").Value = -3240.6667

$ws.Range("H31").Value = 620400.3
$ws.Range("I31").Value = 4871.32
$ws.Range("J31").Value = 1060063.9
$ws.Range("K31").Value = 4871.32
$ws.Range("L31").Value = 1060063.9
$ws.Range("M31").Value = -4576.32
$ws.Range("N31").Value = -1060653.9

$ws.Range("H34").Value = 620400.3
$ws.Range("I34").Value = 4871.32
$ws.Range("J34").Value = 1060063.9
$ws.Range("K34").Value = 4871.32
$ws.Range("L34").Value = 1060063.9
$ws.Range("M34").Value = -4669.32
$ws.Range("N34").Value = -1060467.9

$ws.Range("H113").Value = 2212.2222
$ws.Range("I113").Value = 1985
$ws.Range("J113").Value = 2666.6667
$ws.Range("K113").Value = 1985
$ws.Range("L113").Value = 2666.6667
$ws.Range("M113").Value = 185
$ws.Range("N113").Value = -7006.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2575.582
$ws.Range("I68").Value = 833.125
$ws.Range("J68").Value = 4168.6855
$ws.Range("K68").Value = 2499.375
$ws.Range("L68").Value = 12506.0565
$ws.Range("M68").Value = -1688.375
$ws.Range("N68").Value = -14128.0565

$ws.Range("H71").Value = 2575.582
$ws.Range("I71").Value = 833.125
$ws.Range("J71").Value = 4168.6855
$ws.Range("K71").Value = 7498.125
$ws.Range("L71").Value = 37518.1695
$ws.Range("M71").Value = -3442.125
$ws.Range("N71").Value = -45630.1695

$ws.Range("H86").Value = 782.3333
$ws.Range("I86").Value = 964.6667
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 2894.0001
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -1708.0001
$ws.Range("N86").Value = -4172

$ws.Range("H89").Value = 782.3333
$ws.Range("I89").Value = 964.6667
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 8682.0003
$ws.Range("L89").Value = 5400
$ws.Range("M89").Value = -2754.0003
$ws.Range("N89").Value = -17256

$ws.Range("H109").Value = 3410.125
$ws.Range("J109").Value = 3806
$ws.Range("L109").Value = 11418
$ws.Range("N109").Value = -13498

$ws.Range("H113").Value = 533.2532
$ws.Range("I113").Value = 532.6724
$ws.Range("J113").Value = 534.8570999999999
$ws.Range("K113").Value = 1598.0172
$ws.Range("L113").Value = 1604.5713
$ws.Range("M113").Value = 571.9827999999998
$ws.Range("N113").Value = -5944.5713

$ws.Range("H134").Value = 4687.933
$ws.Range("I134").Value = 6217.273
$ws.Range("J134").Value = 3802.5264
$ws.Range("K134").Value = 18651.819
$ws.Range("L134").Value = 11407.5792
$ws.Range("M134").Value = -13581.819
$ws.Range("N134").Value = -21547.5792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H123").Value = 10325.833
$ws.Range("J123").Value = 10325.833
$ws.Range("L123").Value = 10325.833
$ws.Range("N123").Value = -15225.833

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4875.5
$ws.Range("J22").Value = 6334
$ws.Range("L22").Value = 6334
$ws.Range("N22").Value = -6924

$ws.Range("H27").Value = 4875.5
$ws.Range("J27").Value = 6334
$ws.Range("L27").Value = 6334
$ws.Range("N27").Value = -6548

$ws.Range("H132").Value = 2171.9707
$ws.Range("I132").Value = 1546.619
$ws.Range("J132").Value = 3182.1538
$ws.Range("K132").Value = 4639.857
$ws.Range("L132").Value = 9546.4614
$ws.Range("M132").Value = -2109.857
$ws.Range("N132").Value = -14606.4614

$ws.Range("H136").Value = 5974.84
$ws.Range("I136").Value = 6287.1577
$ws.Range("J136").Value = 4985.8335
$ws.Range("K136").Value = 18861.4731
$ws.Range("L136").Value = 14957.5005
$ws.Range("M136").Value = -16311.4731
$ws.Range("N136").Value = -20057.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1970.0769
$ws.Range("I132").Value = 2286.5173
$ws.Range("J132").Value = 1052.4
$ws.Range("K132").Value = 6859.5519
$ws.Range("L132").Value = 3157.2
$ws.Range("M132").Value = -4329.5519
$ws.Range("N132").Value = -8217.200000000001
